$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 259, pushing the existing 259:265 block down to 262:268.
$ws.Rows("259:261").Insert()

# --- Row 259: new weekly entry (Especial) ---
$ws.Range("A259").Value = 6
$ws.Range("B259").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C259").Value = "Metropolitana"
$ws.Range("D259").Value = 44714
$ws.Range("E259").Value = 13
$ws.Range("F259").Value = 100112043
$ws.Range("G259").Value = "Pepino dulce"
$ws.Range("H259").Value = "Cultivar IV Región"
$ws.Range("I259").Value = "Especial"
$ws.Range("J259").Value = 240
$ws.Range("K259").Value = 16000
$ws.Range("L259").Value = 16000
$ws.Range("M259").Value = 16000
$ws.Range("N259").Value = "$/bandeja 18 kilos"
$ws.Range("O259").Value = "Provincia de Limarí"
$ws.Range("P259").Value = 889
$ws.Range("Q259").Value = 18
$ws.Range("R259").Value = "Hortaliza"

# --- Row 260: new weekly entry (Primera) ---
$ws.Range("A260").Value = 6
$ws.Range("B260").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C260").Value = "Metropolitana"
$ws.Range("D260").Value = 44714
$ws.Range("E260").Value = 13
$ws.Range("F260").Value = 100112043
$ws.Range("G260").Value = "Pepino dulce"
$ws.Range("H260").Value = "Cultivar IV Región"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 350
$ws.Range("K260").Value = 14000
$ws.Range("L260").Value = 14000
$ws.Range("M260").Value = 14000
$ws.Range("N260").Value = "$/bandeja 18 kilos"
$ws.Range("O260").Value = "Provincia de Limarí"
$ws.Range("P260").Value = 778
$ws.Range("Q260").Value = 18
$ws.Range("R260").Value = "Hortaliza"

# --- Row 261: new weekly entry (Segunda) ---
$ws.Range("A261").Value = 6
$ws.Range("B261").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C261").Value = "Metropolitana"
$ws.Range("D261").Value = 44714
$ws.Range("E261").Value = 13
$ws.Range("F261").Value = 100112043
$ws.Range("G261").Value = "Pepino dulce"
$ws.Range("H261").Value = "Cultivar IV Región"
$ws.Range("I261").Value = "Segunda"
$ws.Range("J261").Value = 270
$ws.Range("K261").Value = 11000
$ws.Range("L261").Value = 11000
$ws.Range("M261").Value = 11000
$ws.Range("N261").Value = "$/bandeja 18 kilos"
$ws.Range("O261").Value = "Provincia de Limarí"
$ws.Range("P261").Value = 611
$ws.Range("Q261").Value = 18
$ws.Range("R261").Value = "Hortaliza"
